$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A18").Value = "release/8.0.15"
$ws.Range("B18").Value = "X"
$ws.Range("C18").Value = "X"
$ws.Range("D18").Value = "X"
$ws.Range("E18").Value = "X"
